# ITC and ptc calculations for advanced nuclear
# Update the "Run status" sheet (second sheet) to insert a new "Baseline"
# row at the top of the run-matrix, shift the existing run labels down by
# one row, and clear out the per-state "x"/"o" marker columns (B:F) for
# all rows except the new Baseline row (which only keeps an "r" marker in
# column B). Also clears/repositions the small status legend in H:I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run status")

# --- Column A: run labels (row 2 gets a new "Baseline" entry, and all
#     subsequent labels shift down by one row; a new row 18 is added for
#     the label that used to be on row 17). ---
$labels = @(
  "Baseline",
  "Ref",
  "SMR 20",
  "SMR 100",
  "CO2 low",
  "CO2 high",
  "PTC 000",
  "PTC 100",
  "PTC 270",
  "CAPEX 0.75",
  "CAPEX 1.25",
  "Syn. 0.75",
  "Syn. 1.25",
  "Elec 0.75",
  "Elec 1.25",
  "OM 0.75",
  "OM 1.25"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# --- Columns B:F: clear the per-state marker grid for every row, then
#     re-set column B for the Baseline and Ref rows to "r". ---
$ws.Range("B2:F18").ClearContents()
$ws.Range("B2").Value = "r"
$ws.Range("B3").Value = "r"

# --- Columns H:I: small status legend shifts down one row, and gains a
#     new "o" / "sweep.csv +cashflows" entry on row 3. ---
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

$ws.Range("H3").Value = "o"
$ws.Range("I3").Value = "sweep.csv +cashflows"
$ws.Range("H4").Value = "r"
$ws.Range("I4").Value = "running"

# --- Update the used range / selection bookkeeping to match row 18 being
#     the new last row, with the active selection on F16. ---
$ws.Range("F16").Select()

$wb.Save()
